$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 20,10
$arr[0,0] = -15.9193289081869
$arr[0,1] = 0.07637200578347599
$arr[0,2] = -15.9193289081869
$arr[0,3] = -15.9193289081869
$arr[0,4] = -15.9193289081869
$arr[0,5] = -15.9193289081869
$arr[0,6] = -15.9193289081869
$arr[0,7] = -15.9193289081869
$arr[0,8] = -15.9193289081869
$arr[0,9] = -15.9193289081869
$arr[1,0] = -15.9193289081869
$arr[1,1] = -15.9193289081869
$arr[1,2] = -15.9193289081869
$arr[1,3] = -15.9193289081869
$arr[1,4] = -15.9193289081869
$arr[1,5] = -15.9193289081869
$arr[1,6] = -15.9193289081869
$arr[1,7] = -15.9193289081869
$arr[1,8] = -15.9193289081869
$arr[1,9] = -15.9193289081869
$arr[2,0] = -15.9193289081869
$arr[2,1] = -0.01758991821785622
$arr[2,2] = 0.3788231806228136
$arr[2,3] = -15.9193289081869
$arr[2,4] = 3.910519227418201
$arr[2,5] = -15.9193289081869
$arr[2,6] = 0.6781562918772041
$arr[2,7] = -15.9193289081869
$arr[2,8] = 1.965205339857661
$arr[2,9] = -15.9193289081869
$arr[3,0] = -15.9193289081869
$arr[3,1] = 0.6081022798491682
$arr[3,2] = -15.9193289081869
$arr[3,3] = -15.9193289081869
$arr[3,4] = -15.9193289081869
$arr[3,5] = 3.353339600589717
$arr[3,6] = -15.9193289081869
$arr[3,7] = -15.9193289081869
$arr[3,8] = -15.9193289081869
$arr[3,9] = -15.9193289081869
$arr[4,0] = -15.9193289081869
$arr[4,1] = -15.9193289081869
$arr[4,2] = -15.9193289081869
$arr[4,3] = -15.9193289081869
$arr[4,4] = -15.9193289081869
$arr[4,5] = -15.9193289081869
$arr[4,6] = -15.9193289081869
$arr[4,7] = -15.9193289081869
$arr[4,8] = -15.9193289081869
$arr[4,9] = -15.9193289081869
$arr[5,0] = 3.088828854855712
$arr[5,1] = -15.9193289081869
$arr[5,2] = -15.9193289081869
$arr[5,3] = -15.9193289081869
$arr[5,4] = -15.9193289081869
$arr[5,5] = -15.9193289081869
$arr[5,6] = -15.9193289081869
$arr[5,7] = -15.9193289081869
$arr[5,8] = -15.9193289081869
$arr[5,9] = -15.9193289081869
$arr[6,0] = -15.9193289081869
$arr[6,1] = -15.9193289081869
$arr[6,2] = -15.9193289081869
$arr[6,3] = 1.823487077158004
$arr[6,4] = -15.9193289081869
$arr[6,5] = -15.9193289081869
$arr[6,6] = -15.9193289081869
$arr[6,7] = -15.9193289081869
$arr[6,8] = -15.9193289081869
$arr[6,9] = -15.9193289081869
$arr[7,0] = 3.522514997953153
$arr[7,1] = -15.9193289081869
$arr[7,2] = -15.9193289081869
$arr[7,3] = -15.9193289081869
$arr[7,4] = -15.9193289081869
$arr[7,5] = -15.9193289081869
$arr[7,6] = -15.9193289081869
$arr[7,7] = -15.9193289081869
$arr[7,8] = -15.9193289081869
$arr[7,9] = -15.9193289081869
$arr[8,0] = -15.9193289081869
$arr[8,1] = -15.9193289081869
$arr[8,2] = -15.9193289081869
$arr[8,3] = -15.9193289081869
$arr[8,4] = -15.9193289081869
$arr[8,5] = -15.9193289081869
$arr[8,6] = -15.9193289081869
$arr[8,7] = 4.321905978918958
$arr[8,8] = -15.9193289081869
$arr[8,9] = 1.775335040396872
$arr[9,0] = -15.9193289081869
$arr[9,1] = -15.9193289081869
$arr[9,2] = -15.9193289081869
$arr[9,3] = 2.213276248393478
$arr[9,4] = -15.9193289081869
$arr[9,5] = 2.30996402791783
$arr[9,6] = -15.9193289081869
$arr[9,7] = -15.9193289081869
$arr[9,8] = -15.9193289081869
$arr[9,9] = 2.252888595848371
$arr[10,0] = -15.9193289081869
$arr[10,1] = -15.9193289081869
$arr[10,2] = -15.9193289081869
$arr[10,3] = -15.9193289081869
$arr[10,4] = -15.9193289081869
$arr[10,5] = -15.9193289081869
$arr[10,6] = -15.9193289081869
$arr[10,7] = -15.9193289081869
$arr[10,8] = -15.9193289081869
$arr[10,9] = -15.9193289081869
$arr[11,0] = -15.9193289081869
$arr[11,1] = -15.9193289081869
$arr[11,2] = -15.9193289081869
$arr[11,3] = 1.779274030673821
$arr[11,4] = -15.9193289081869
$arr[11,5] = -15.9193289081869
$arr[11,6] = -15.9193289081869
$arr[11,7] = -15.9193289081869
$arr[11,8] = 2.195699497004713
$arr[11,9] = 2.017659771011587
$arr[12,0] = -15.9193289081869
$arr[12,1] = -15.9193289081869
$arr[12,2] = 1.525596995481486
$arr[12,3] = -15.9193289081869
$arr[12,4] = -15.9193289081869
$arr[12,5] = -15.9193289081869
$arr[12,6] = -15.9193289081869
$arr[12,7] = -15.9193289081869
$arr[12,8] = -15.9193289081869
$arr[12,9] = 1.513906052676234
$arr[13,0] = -15.9193289081869
$arr[13,1] = -15.9193289081869
$arr[13,2] = 0.4673587895298714
$arr[13,3] = -15.9193289081869
$arr[13,4] = -15.9193289081869
$arr[13,5] = -15.9193289081869
$arr[13,6] = -15.9193289081869
$arr[13,7] = -15.9193289081869
$arr[13,8] = -15.9193289081869
$arr[13,9] = -15.9193289081869
$arr[14,0] = -15.9193289081869
$arr[14,1] = -15.9193289081869
$arr[14,2] = -15.9193289081869
$arr[14,3] = -15.9193289081869
$arr[14,4] = -15.9193289081869
$arr[14,5] = -15.9193289081869
$arr[14,6] = -15.9193289081869
$arr[14,7] = -15.9193289081869
$arr[14,8] = 2.549251386147876
$arr[14,9] = -15.9193289081869
$arr[15,0] = -15.9193289081869
$arr[15,1] = 0.4312132290429203
$arr[15,2] = 0.1104112470164083
$arr[15,3] = -15.9193289081869
$arr[15,4] = -15.9193289081869
$arr[15,5] = -15.9193289081869
$arr[15,6] = 2.209657206671276
$arr[15,7] = -15.9193289081869
$arr[15,8] = 1.707813072223102
$arr[15,9] = -15.9193289081869
$arr[16,0] = -15.9193289081869
$arr[16,1] = -15.9193289081869
$arr[16,2] = -15.9193289081869
$arr[16,3] = -15.9193289081869
$arr[16,4] = -15.9193289081869
$arr[16,5] = -15.9193289081869
$arr[16,6] = 1.58422511901383
$arr[16,7] = -15.9193289081869
$arr[16,8] = 1.259423463456497
$arr[16,9] = -15.9193289081869
$arr[17,0] = -15.9193289081869
$arr[17,1] = -15.9193289081869
$arr[17,2] = 2.74071226280173
$arr[17,3] = -15.9193289081869
$arr[17,4] = -15.9193289081869
$arr[17,5] = -15.9193289081869
$arr[17,6] = 2.238259811509393
$arr[17,7] = -15.9193289081869
$arr[17,8] = -15.9193289081869
$arr[17,9] = -15.9193289081869
$arr[18,0] = -15.9193289081869
$arr[18,1] = 3.116364262747852
$arr[18,2] = 2.738621001000316
$arr[18,3] = -15.9193289081869
$arr[18,4] = 2.31090266838535
$arr[18,5] = -15.9193289081869
$arr[18,6] = 1.988858449845344
$arr[18,7] = -15.9193289081869
$arr[18,8] = -15.9193289081869
$arr[18,9] = 2.294292718660835
$arr[19,0] = -15.9193289081869
$arr[19,1] = 2.680979673853594
$arr[19,2] = -15.9193289081869
$arr[19,3] = 3.068770474683313
$arr[19,4] = -15.9193289081869
$arr[19,5] = 2.26930339851751
$arr[19,6] = 1.062160332648851
$arr[19,7] = -15.9193289081869
$arr[19,8] = -15.9193289081869
$arr[19,9] = -15.9193289081869

$ws.Range("B2:K21").Value = $arr
